# Fix circular-reference SUM ranges on the LOE sheet's summary rows
# (Management / Project Management / TOTAL rows). Each formula previously
# summed a range that looped back on the formula's own row (e.g. E17:E43
# includes rows 30/31/32), causing a circular reference. Point the SUMs at
# the actual task-data rows (3:29/3:30/3:31) instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LOE")

$ws.Range("D30").Formula = "=ROUND(SUM(E3:E29)*`$C`$30,0)"
$ws.Range("E30").Formula = "=ROUND(SUM(E3:E29)*0.25,0)"

$ws.Range("D31").Formula = "=ROUND(SUM(E3:E30)*0.20,0)"
$ws.Range("E31").Formula = "=ROUND(SUM(E3:E30)*0.20,0)"

$ws.Range("E32").Formula = "=SUM(E3:E31)"
$ws.Range("G32").Formula = "=TEXT(SUM(G3:G31),`"`$#,##0`")"
